$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value  = 254
$ws.Range("I2").Value  = 727
$ws.Range("J2").Value  = 3053
$ws.Range("K2").Value  = 17
$ws.Range("L2").Value  = 767
$ws.Range("M2").Value  = 58
$ws.Range("N2").Value  = 515
$ws.Range("O2").Value  = 0
$ws.Range("P2").Value  = 15
$ws.Range("Q2").Value  = 4
$ws.Range("R2").Value  = 27
$ws.Range("S2").Value  = 315
$ws.Range("T2").Value  = 516
$ws.Range("U2").Value  = 48
$ws.Range("V2").Value  = 4644
$ws.Range("W2").Value  = 1
$ws.Range("X2").Value  = 4529
$ws.Range("Y2").Value  = 4
$ws.Range("Z2").Value  = 76
$ws.Range("AA2").Value = 43
